$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1 "Play Combat Masters Free - Slot Game Review").
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Build the exact run layout (leading empty run + bold "Meta description" run
# + plain run with the rest of the sentence) via InsertXML so the paragraph
# mirrors the structure used elsewhere in the document.
$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Experience the engaging gameplay system of Combat Masters! Read our review and play for free to enjoy medium volatility and fascinating graphics.</w:t></w:r></w:p>' + `
  '<w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated "Play Combat Masters Free - Slot Game Review"
#    paragraph near the end of the document (its content now lives in the
#    new Meta description paragraph above). Search for it starting after the
#    title (so the first, legitimate, occurrence is skipped), then delete
#    the whole Paragraph object (not a loose Range) so the paragraph mark is
#    removed along with the text and no blank paragraph is left behind.
# ---------------------------------------------------------------------------
$searchRange = $d.Range($title.Range.End, $d.Content.End)
$found = $searchRange.Find.Execute("Play Combat Masters Free - Slot Game Review", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $total = $d.Paragraphs.Count
    for ($i = 1; $i -le $total; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Start -le $searchRange.Start -and $candidate.Range.End -ge $searchRange.End) {
            $candidate.Range.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the italic closing paragraph's text with the new feature-image
#    prompt, keeping its existing (italic) run formatting intact. Using an
#    explicit Document.Range(start, end) (rather than Find/Replace or the
#    Paragraph's own Range object) performs a clean replace and avoids
#    smart-quote autocorrection of the straight apostrophes in the new text.
# ---------------------------------------------------------------------------
$total2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($total2)
for ($i = 1; $i -le $total2; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Experience the engaging gameplay system of Combat Masters*") {
        $lastPara = $candidate
        break
    }
}

$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End
$lastRange = $d.Range($lastStart, $lastEnd)
$lastRange.Text = "Create a feature image for Combat Masters: Design a cartoon-style image featuring a happy and fierce-looking Maya warrior wearing glasses. The Maya warrior should be holding a sword in one hand and a shield in the other, ready for battle. The shield should have the game's logo emblazoned on it. The background should be an epic battlefield with medieval castles, mountains, and a fantasy landscape. The colors should be bright and vibrant, with a focus on blue and purple hues. Additionally, include the game's title, Combat Masters, in a bold, stylized font."

Write-Output "done"
